$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3679.138
$ws.Range("I64").Value = 3659.92
$ws.Range("K64").Value = 3659.92
$ws.Range("M64").Value = -3411.92
$ws.Range("H67").Value = 3679.138
$ws.Range("I67").Value = 3659.92
$ws.Range("K67").Value = 3659.92
$ws.Range("M67").Value = -2801.92
$ws.Range("H100").Value = 2445.75
$ws.Range("J100").Value = 3424.25
$ws.Range("L100").Value = 3424.25
$ws.Range("N100").Value = -4506.25
$ws.Range("H107").Value = 341.33334
$ws.Range("I107").Value = 294.2857
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 294.2857
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1625.7143
$ws.Range("N107").Value = -4840
$ws.Range("H132").Value = 1868.6552
$ws.Range("I132").Value = 1968.8846
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 5906.6538
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -3376.6538
$ws.Range("N132").Value = -8060
$ws.Range("H135").Value = 2078.6843
$ws.Range("I135").Value = 1838.3846
$ws.Range("J135").Value = 2599.3333
$ws.Range("K135").Value = 16545.4614
$ws.Range("L135").Value = 23393.9997
$ws.Range("M135").Value = -14010.4614
$ws.Range("N135").Value = -28463.9997
$ws.Range("H138").Value = 1436.7715
$ws.Range("I138").Value = 1264.3226
$ws.Range("J138").Value = 2773.25
$ws.Range("K138").Value = 3792.9678
$ws.Range("L138").Value = 8319.75
$ws.Range("M138").Value = 1347.0322
$ws.Range("N138").Value = -18599.75
$ws.Range("H141").Value = 3338.5
$ws.Range("I141").Value = 2895.25
$ws.Range("K141").Value = 8685.75
$ws.Range("M141").Value = -3505.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 7357478
$ws.Range("I45").Value = 2717.7778
$ws.Range("J45").Value = 15631584
$ws.Range("K45").Value = 2717.7778
$ws.Range("L45").Value = 15631584
$ws.Range("M45").Value = -2340.7778
$ws.Range("N45").Value = -15632338
$ws.Range("H74").Value = 38150.535
$ws.Range("I74").Value = 72888.14
$ws.Range("J74").Value = 3412.9285
$ws.Range("K74").Value = 72888.14
$ws.Range("L74").Value = 3412.9285
$ws.Range("M74").Value = -72014.14
$ws.Range("N74").Value = -5160.9285
$ws.Range("H77").Value = 38150.535
$ws.Range("I77").Value = 72888.14
$ws.Range("J77").Value = 3412.9285
$ws.Range("K77").Value = 364440.7
$ws.Range("L77").Value = 17064.6425
$ws.Range("M77").Value = -360072.7
$ws.Range("N77").Value = -25800.6425
$ws.Range("H94").Value = 34317
$ws.Range("J94").Value = 34317
$ws.Range("L94").Value = 34317
$ws.Range("N94").Value = -36119
$ws.Range("H96").Value = 34999.668
$ws.Range("J96").Value = 34999.668
$ws.Range("L96").Value = 34999.668
$ws.Range("N96").Value = -40491.668
$ws.Range("H102").Value = 158647.28
$ws.Range("I102").Value = 168420
$ws.Range("K102").Value = 168420
$ws.Range("M102").Value = -166798
$ws.Range("H132").Value = 1834.2222
$ws.Range("I132").Value = 1251.3
$ws.Range("J132").Value = 3499.7144
$ws.Range("K132").Value = 3753.9
$ws.Range("L132").Value = 10499.1432
$ws.Range("M132").Value = -1223.9
$ws.Range("N132").Value = -15559.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2330.7
$ws.Range("I107").Value = 1975
$ws.Range("K107").Value = 1975
$ws.Range("M107").Value = -55
$ws.Range("H135").Value = 62501.145
$ws.Range("J135").Value = 62501.145
$ws.Range("L135").Value = 62501.145
$ws.Range("N135").Value = -72641.14499999999
$ws.Range("H140").Value = 43292.715
$ws.Range("J140").Value = 43471.8
$ws.Range("L140").Value = 43471.8
$ws.Range("N140").Value = -53831.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1743.5428
$ws.Range("I58").Value = 1576.1538
$ws.Range("J58").Value = 2227.111
$ws.Range("K58").Value = 1576.1538
$ws.Range("L58").Value = 2227.111
$ws.Range("M58").Value = -1373.1538
$ws.Range("N58").Value = -2633.111
$ws.Range("H86").Value = 5060.2085
$ws.Range("I86").Value = 4227.9287
$ws.Range("J86").Value = 6225.4
$ws.Range("K86").Value = 4227.9287
$ws.Range("L86").Value = 6225.4
$ws.Range("M86").Value = -3104.9287
$ws.Range("N86").Value = -8471.4
$ws.Range("H89").Value = 5060.2085
$ws.Range("I89").Value = 4227.9287
$ws.Range("J89").Value = 6225.4
$ws.Range("K89").Value = 21139.6435
$ws.Range("L89").Value = 31127
$ws.Range("M89").Value = -15523.6435
$ws.Range("N89").Value = -42359
$ws.Range("H132").Value = 2347.8572
$ws.Range("I132").Value = 2347.8572
$ws.Range("K132").Value = 7043.571599999999
$ws.Range("M132").Value = -4513.571599999999
$ws.Range("H134").Value = 57191.668
$ws.Range("I134").Value = 1268.6364
$ws.Range("J134").Value = 145070.72
$ws.Range("K134").Value = 3805.9092
$ws.Range("L134").Value = 435212.16
$ws.Range("M134").Value = -1270.9092
$ws.Range("N134").Value = -440282.16
$ws.Range("H136").Value = 1743.5428
$ws.Range("I136").Value = 1576.1538
$ws.Range("J136").Value = 2227.111
$ws.Range("K136").Value = 4728.4614
$ws.Range("L136").Value = 6681.333
$ws.Range("M136").Value = -2178.4614
$ws.Range("N136").Value = -11781.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 37.1875
$ws.Range("I2").Value = 37.1875
$ws.Range("K2").Value = 223.125
$ws.Range("M2").Value = -110.125
$ws.Range("H38").Value = 84
$ws.Range("I38").Value = 23
$ws.Range("K38").Value = 69
$ws.Range("M38").Value = 278
$ws.Range("H39").Value = 3295.3333
$ws.Range("J39").Value = 3465.647
$ws.Range("L39").Value = 10396.941
$ws.Range("N39").Value = -10984.941
$ws.Range("H113").Value = 42538.707
$ws.Range("I113").Value = 733.9
$ws.Range("J113").Value = 72399.28999999999
$ws.Range("K113").Value = 2201.7
$ws.Range("L113").Value = 217197.87
$ws.Range("M113").Value = -31.69999999999982
$ws.Range("N113").Value = -221537.87
$ws.Range("H122").Value = 360.13333
$ws.Range("I122").Value = 510.75
$ws.Range("J122").Value = 305.36365
$ws.Range("K122").Value = 4596.75
$ws.Range("L122").Value = 2748.27285
$ws.Range("M122").Value = -2146.75
$ws.Range("N122").Value = -7648.27285
$ws.Range("H140").Value = 2202.647
$ws.Range("I140").Value = 1235.9546
$ws.Range("K140").Value = 3707.8638
$ws.Range("M140").Value = 1472.1362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7309.4
$ws.Range("I70").Value = 7306.778
$ws.Range("J70").Value = 7313.3335
$ws.Range("K70").Value = 7306.778
$ws.Range("L70").Value = 7313.3335
$ws.Range("M70").Value = -7036.778
$ws.Range("N70").Value = -7853.3335
$ws.Range("H73").Value = 7309.4
$ws.Range("I73").Value = 7306.778
$ws.Range("J73").Value = 7313.3335
$ws.Range("K73").Value = 7306.778
$ws.Range("L73").Value = 7313.3335
$ws.Range("M73").Value = -6370.778
$ws.Range("N73").Value = -9185.333500000001
$ws.Range("H102").Value = 2210.7646
$ws.Range("I102").Value = 2265.5833
$ws.Range("J102").Value = 2079.2
$ws.Range("K102").Value = 2265.5833
$ws.Range("L102").Value = 2079.2
$ws.Range("M102").Value = -643.5832999999998
$ws.Range("N102").Value = -5323.2
$ws.Range("H122").Value = 88247.19500000001
$ws.Range("I122").Value = 98148.35000000001
$ws.Range("K122").Value = 294445.05
$ws.Range("M122").Value = -291995.05
$ws.Range("H132").Value = 5706.5264
$ws.Range("I132").Value = 4342.231
$ws.Range("K132").Value = 13026.693
$ws.Range("M132").Value = -10496.693
$ws.Range("H135").Value = 70000
$ws.Range("J135").Value = 70000
$ws.Range("L135").Value = 70000
$ws.Range("N135").Value = -80140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4611.115
$ws.Range("J22").Value = 9153.25
$ws.Range("L22").Value = 9153.25
$ws.Range("N22").Value = -9743.25
$ws.Range("H27").Value = 4611.115
$ws.Range("J27").Value = 9153.25
$ws.Range("L27").Value = 9153.25
$ws.Range("N27").Value = -9367.25
$ws.Range("H40").Value = 4634608
$ws.Range("I40").Value = 5156.3335
$ws.Range("J40").Value = 18522962
$ws.Range("K40").Value = 5156.3335
$ws.Range("L40").Value = 18522962
$ws.Range("M40").Value = -5020.3335
$ws.Range("N40").Value = -18523234
$ws.Range("H136").Value = 1862.5758
$ws.Range("I136").Value = 1505.9131
$ws.Range("J136").Value = 2682.9
$ws.Range("K136").Value = 4517.7393
$ws.Range("L136").Value = 8048.700000000001
$ws.Range("M136").Value = -1967.7393
$ws.Range("N136").Value = -13148.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4397481
$ws.Range("I96").Value = 13900
$ws.Range("K96").Value = 13900
$ws.Range("M96").Value = -12527
$ws.Range("H126").Value = 112590.5
$ws.Range("I126").Value = 220631
$ws.Range("J126").Value = 4550
$ws.Range("K126").Value = 661893
$ws.Range("L126").Value = 13650
$ws.Range("M126").Value = -659423
$ws.Range("N126").Value = -18590
$ws.Range("H132").Value = 15529
$ws.Range("I132").Value = 19897.766
$ws.Range("J132").Value = 3150.8333
$ws.Range("K132").Value = 59693.298
$ws.Range("L132").Value = 9452.499899999999
$ws.Range("M132").Value = -57163.298
$ws.Range("N132").Value = -14512.4999
